$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows 3 through 11 (old extra data rows) and column E entirely
$ws.Rows("3:11").Delete()
$ws.Columns("E:E").Delete()

# Update header row
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("C1").Value = "Distancia"
$ws.Range("D1").Value = "Tempo"

# Update data row
$ws.Range("A2").Value = 81
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 3231
$ws.Range("D2").Value = 0.04646468162536621
